# ENACOM_GH.xlsx update - refreshed figures for the "DIRECTV ARGENTINA S.A."
# row (renamed from "DIRECTV ARGENTINA S.A./AT&T") plus several other
# updated revenue figures in column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Rename the DIRECTV group label (row 5)
$ws.Range("A5").Value = "DIRECTV ARGENTINA S.A."

# Updated "Ingresos s/TCFV" (column B) figures
$ws.Range("B5").Value = 1785083970.48
$ws.Range("B10").Value = 11908070231.99
$ws.Range("B11").Value = 3951631352.2000003
$ws.Range("B12").Value = 4122487164.6999998
$ws.Range("B13").Value = 6804835332.3099995
$ws.Range("B14").Value = 6068088682.04
$ws.Range("B16").Value = 3554471505.0800004

# Match the author's last on-screen selection before saving
$ws.Range("C20").Select()
